# RPA datasets push 2024-06-15
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 01_IB전략컨설팅부 -------------------------------------------
# Remove the 신한제12호스팩 row (row 14) - it has dropped out of this push.
$ws1 = $wb.Worksheets.Item("01_IB전략컨설팅부")
$ws1.Rows.Item(14).Delete()

# --- Sheet 2: 02_38커뮤니케이션(최근일자기준) ------------------------------
$ws2 = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

# Insert a new deal (뱅크웨어글로벌) right after row 2 (아이빔테크놀로지),
# pushing everything else down by one row.
$ws2.Rows.Item(3).Insert()
$ws2.Cells.Item(3, 1).Value = "뱅크웨어글로벌"
$ws2.Cells.Item(3, 2).Value = "2024.07.08~07.12"
$ws2.Cells.Item(3, 3).Value = "16,000~19,000"
$ws2.Cells.Item(3, 4).Value = "-"
$ws2.Cells.Item(3, 5).Value = 22400
$ws2.Cells.Item(3, 6).Value = "미래에셋증권"

# 시프트업(유가) has completed its listing and drops off the bottom of
# the tracked list (was row 21, now row 22 after the insert above).
$ws2.Rows.Item(22).Delete()

# 한국스팩14호's underwriter info was updated this push.
$ws2.Cells.Item(21, 6).Value = "KB증권"
